$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 2
$ws.Range("F3").Value = -3
$ws.Range("H3").Value = 46

$ws.Range("C3").Select()
